# Append a new run's worth of scraped job listings to the "ランサーズ" sheet.
#  - Insert a new row for a brand-new listing right after the header (row 2)
#  - Re-stamp every existing row with the new scrape timestamp
#  - Append two more brand-new rows at the bottom
#  - Widen columns D and H slightly
#  - Rebuild the hyperlinks on column F so they keep following their URL

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-12-28 18:26:33"

# Final contents for data rows 2..8 (row 1 is the header and is unchanged).
$rows = @(
    @{ B = "Bubble × AI API組込みのWebアプリ開発パートナー募集";
       C = "システム開発";
       D = "1,000,000 円 ~ 3,000,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5462677";
       G = 545;
       H = "🔥AI,API ◆開発 ◇アプリ" },
    @{ B = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪";
       C = "システム開発";
       D = "20,000 円 ~ 50,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5217096";
       G = 243;
       H = "🔥API ◆ツール" },
    @{ B = "GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集";
       C = "システム開発";
       D = "500,000 円 ~ 1,000,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5457458";
       G = 75;
       H = "◆開発" },
    @{ B = "【急募】WordoressサイトスピードUPのための専門家を探しています!";
       C = "システム開発";
       D = "20,000 円 ~ 50,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5462581";
       G = 33;
       H = "◇サイト" },
    @{ B = "FXレイテンシーアービトラージの検証(環境設計・比較評価・PoC)";
       C = "システム開発";
       D = "300,000 円 ~ 500,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5462397";
       G = 25;
       H = "" },
    @{ B = "【急募】簡単なHP作成とAWS構築をしてくれる方募集(オンラインMTG参加必須)";
       C = "システム開発";
       D = "50,000 円 ~ 100,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5462522";
       G = 18;
       H = "" },
    @{ B = "【急募】VPS(Ubuntu)フロント・バック統合/Docker構築・移行";
       C = "システム開発";
       D = "10,000 円 ~ 20,000 円 / 固定";
       E = "期限情報なし";
       F = "https://www.lancers.jp/work/detail/5462712";
       G = 10;
       H = "" }
)

# Remove the existing hyperlink relationships; we rebuild them from scratch
# below once all the URLs are back in the right rows, so each one keeps
# pointing at the correct target after the shift.
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    if ($data.H -ne "") {
        $ws.Cells.Item($r, 8).Value = $data.H
    }
}

# Re-create hyperlinks for column F, row by row, in the same top-to-bottom
# order as the data so relationship ids line up the same way Excel would
# assign them.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $rows[$i].F)
}

# Adding a hyperlink resets the cell style to a fresh copy of "Hyperlink";
# reapply the named style explicitly so column F keeps reusing the original
# style record instead of accumulating duplicates.
$ws.Range("F2:F8").Style = "Hyperlink"

# Column width tweaks (raw OOXML width units 30->32 and 12->18). Excel's
# ColumnWidth property is expressed in characters and differs from the
# stored width by the fixed padding of 5/6 of a character, so compensate
# for that offset to land exactly on the target raw widths.
$ws.Columns.Item(4).ColumnWidth = 32 - 5/6
$ws.Columns.Item(8).ColumnWidth = 18 - 5/6
